$wb = $excel.ActiveWorkbook

$wsCont = $wb.Worksheets.Item("Cont adminstrator")
$wsCont.Range("B15").Value = "Scoala particulara Miranda D"
$wsCont.Range("C15").Value = "mirandascoalaa1@automation.33mail.com"

$wsRec = $wb.Worksheets.Item("Receptie")
$wsRec.Range("B2").Value = "artoise23@staffcalendis.33mail.com"
$wsRec.Range("B3").Value = "artemisa5@staffcalendis.33mail.com"
$wsRec.Range("B4").Value = "goblins1@staffcalendis.33mail.com"

$wsAng = $wb.Worksheets.Item("Angajati")
$wsAng.Range("B2").Value = "fantastricarep@staffcalendis.33mail.com"
$wsAng.Range("B3").Value = "demokratiks4@staffcalendis.33mail.com"
$wsAng.Range("B4").Value = "zanzibara@staffcalendis.33mail.com"
$wsAng.Range("B5").Value = "zebralda@staffcalendis.33mail.com"
